$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency market data (price & 1h volume change)
# Cells whose new text could be misread as a number are forced to text
# via NumberFormat "@" before assignment, so Excel keeps them as strings
# (matching the original inlineStr/shared-string cell type).

$ws.Range("D2").Value = "27.180.17"
$ws.Range("E2").Value = "  +0.74%  "

$ws.Range("D3").Value = "1.561.37"
$ws.Range("E3").Value = "  -0.03%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.01"
$ws.Range("E4").Value = "  +0.17%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "210.22"
$ws.Range("E5").Value = "  +1.30%  "

$ws.Range("E6").Value = "  +0.08%  "

$ws.Range("E7").Value = "  -0.27%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "22.05"
$ws.Range("E8").Value = "  -0.34%  "

$ws.Range("E9").Value = "  +0.05%  "

$ws.Range("E10").Value = "  -1.07%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0871"
$ws.Range("E11").Value = "  +1.80%  "

$ws.Range("D12").Value = "1.780.17"
$ws.Range("E12").Value = "  -0.22%  "

$ws.Range("D13").Value = "1.557.14"
$ws.Range("E13").Value = "  -0.36%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.77"
$ws.Range("E14").Value = "  +0.37%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.517"
$ws.Range("E15").Value = "  -0.61%  "

$ws.Range("D16").Value = "27.130.87"
$ws.Range("E16").Value = "  +0.63%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.74"
$ws.Range("E17").Value = "  -0.26%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.45"
$ws.Range("E18").Value = "  +1.28%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "216.35"
$ws.Range("E19").Value = "  +0.28%  "

$ws.Range("D20").Value = "0.0₃0701"
$ws.Range("E20").Value = "  -0.63%  "

$ws.Range("E21").Value = "  -0.13%  "

$ws.Range("E22").Value = "  +0.33%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.20"
$ws.Range("E23").Value = "  +0.19%  "

$ws.Range("E24").Value = "  +0.45%  "

$ws.Range("E25").Value = "  -0.49%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.62"
$ws.Range("E26").Value = "  +0.10%  "

$ws.Range("B27").Value = "Stellar"
$ws.Range("C27").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.107"
$ws.Range("E27").Value = "  +1.83%  "

$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.01"
$ws.Range("E28").Value = "  -0.56%  "

$ws.Range("E29").Value = "  -0.06%  "

$ws.Range("E30").Value = "  +2.04%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0470"
$ws.Range("E31").Value = "  -0.11%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.23"
$ws.Range("E32").Value = "  +0.15%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.15"
$ws.Range("E33").Value = "  +1.27%  "

$ws.Range("D34").Value = "1.436.46"
$ws.Range("E34").Value = "  +0.87%  "

$ws.Range("E35").Value = "  +4.01%  "

$ws.Range("E36").Value = "  +0.40%  "

$ws.Range("E37").Value = "  -0.47%  "

$ws.Range("E38").Value = "  +0.41%  "

$ws.Range("E39").Value = "  -0.27%  "

$ws.Range("E40").Value = "  +1.32%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.806"
$ws.Range("E41").Value = "  -0.26%  "

$ws.Range("E42").Value = "  -0.18%  "

$ws.Range("E43").Value = "  -0.08%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.997"
$ws.Range("E44").Value = "  -0.90%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "64.14"
$ws.Range("E45").Value = "  -0.66%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.74"
$ws.Range("E46").Value = "  -0.54%  "

$ws.Range("D47").Value = "1.693.45"
$ws.Range("E47").Value = "  -0.24%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "85.35"
$ws.Range("E48").Value = "  -2.09%  "

$ws.Range("E49").Value = "  +0.56%  "

$ws.Range("D50").Value = "0.0₆0100"
$ws.Range("E50").Value = "  +0.03%  "

$ws.Range("E51").Value = "  -1.02%  "
